$d = $word.ActiveDocument

# --- Change 1 ---
# The USER STORY paragraph had its text split into two runs with a
# "_GoBack" bookmark sandwiched between them. Running Find/Replace across
# the whole phrase (spanning both runs and the bookmark) collapses
# everything into a single run - carrying over the first run's
# formatting - and drops the bookmark that used to sit in the middle.
$oldStory = "As a caseworker or parent, it would be beneficial to understand what the quality star ratings represent so that I can choose the best provider for a child."
$newStory = "As a caseworker or parent, it would be beneficial to understand what the quality star ratings represent so that I can choose the best provider for a child."
$d.Content.Find.Execute($oldStory, $true, $false, $false, $false, $false, $true, 1, $false, $newStory, 2) | Out-Null

# --- Change 2 ---
# Add a new run of text to the previously empty paragraph right after
# "ACCEPTANCE CRITERIA", then re-plant the single document "_GoBack"
# bookmark immediately after that new text (Word keeps exactly one
# "_GoBack" bookmark, tracking the most recent edit location).
$acceptanceParagraph = $d.Paragraphs.Item(16)
$insertRange = $acceptanceParagraph.Range
$insertRange.Collapse(1)
$insertStart = $insertRange.Start

$newText = "1. User should get the definitions of all the Quality star rating from search result page"

# Type the new text plus one throwaway trailing placeholder character.
# (Placing a zero-length bookmark range exactly on the paragraph-mark
# boundary is mishandled by this runtime, so we keep the bookmark
# insertion point one character away from that boundary by typing an
# extra placeholder char first, then deleting it again afterwards.)
$insertRange.InsertBefore($newText + "X")

$typedRange = $d.Range($insertStart, $insertStart + $newText.Length)
$typedRange.Font.Color = 0

$bookmarkPos = $insertStart + $newText.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# Remove the placeholder character now that the bookmark is anchored.
$placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholderRange.Delete()
